# Sort the calibration data (rows 2-8, columns A:D) by column A (time) ascending.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 7
$numCols = 4

# Read current data into an array of rows (use Value2 - Value getter is unreliable in this interop)
$data = New-Object 'object[,]' $numRows, $numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r, $c] = $ws.Cells.Item($r + 2, $c + 1).Value2
    }
}

# Build (index, key) pairs and sort by key using -Property (scriptblock closures over
# indexed arrays are unreliable in this interop's Sort-Object implementation).
$pairs = @()
for ($r = 0; $r -lt $numRows; $r++) {
    $pairs += [PSCustomObject]@{ Idx = $r; Key = $data[$r, 0] }
}
$sortedPairs = $pairs | Sort-Object -Property Key

# Write back sorted rows
$destRow = 2
foreach ($p in $sortedPairs) {
    $idx = $p.Idx
    for ($c = 0; $c -lt $numCols; $c++) {
        $ws.Cells.Item($destRow, $c + 1).Value2 = $data[$idx, $c]
    }
    $destRow++
}
